$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 126: new weekly data point (update in place).
# Columns A/B/C/E/F/G/H/I/N/O/Q/R stay the same; only the date (D) and the
# volume/price figures (J/K/L/M/P) change for this entry.
$ws.Cells.Item(126, 4).Value = 44628
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 2500
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = 2750
$ws.Cells.Item(126, 16).Value = 1833

# Row 127: newly appended row holding what used to be row 126's data
# (same market/category/etc.), preserving the old date + price figures.
$ws.Cells.Item(127, 1).Value = 8
$ws.Cells.Item(127, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(127, 3).Value = "Coquimbo"
$ws.Cells.Item(127, 4).Value = 44544
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 5).Value = 4
$ws.Cells.Item(127, 6).Value = 100112040
$ws.Cells.Item(127, 7).Value = "Cilantro"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 3200
$ws.Cells.Item(127, 11).Value = 1500
$ws.Cells.Item(127, 12).Value = 2000
$ws.Cells.Item(127, 13).Value = 1750
$ws.Cells.Item(127, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(127, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(127, 16).Value = 1167
$ws.Cells.Item(127, 17).Value = 1.5
$ws.Cells.Item(127, 18).Value = "Hortaliza"
